$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 951
$ws.Range("B2").Value = 865
$ws.Range("C2").Value = 865
$ws.Range("D2").Value = 865
$ws.Range("E2").Value = 961
$ws.Range("F2").Value = 964
$ws.Range("G2").Value = 949
$ws.Range("H2").Value = 997
